# Practise workbook: add "Sheet2" after Sheet1, clear a stray empty cell on
# Sheet1, and make Sheet2 the active sheet/tab.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Sheet1: drop the now-empty D16 cell (it held no value).
$sheet1.Range("D16").ClearContents()

# Insert the new sheet right after Sheet1 and make it active.
$ws = $wb.Worksheets.Add($null, $sheet1)

# --- Row 2 -----------------------------------------------------------
$ws.Range("A2").Value = "name"
$ws.Range("B2").Value = "sales"
$ws.Range("C2").Value = "QUE?"

# --- Row 3 -------------------------------------------------------------
$ws.Range("A3").Value = "RAHUL"
$ws.Range("B3").Value = 1200
$ws.Range("C3").Value = "TOTAL SALES"
$ws.Range("D3").Formula = "=SUM(B3:B7)"
$ws.Range("F3").Value = "name"
$ws.Range("G3").Value = "sales"
$ws.Range("H3").Value = "region "
$ws.Range("I3").Value = "SEARCH"

# --- Row 4 -------------------------------------------------------------
$ws.Range("A4").Value = "PRIYA"
$ws.Range("B4").Value = 950
$ws.Range("C4").Value = "AVERAGE SALE"
$ws.Range("D4").Formula = "=AVERAGE(B3:B7)"
$ws.Range("F4").Value = "RAHUL"
$ws.Range("G4").Value = 1200
$ws.Range("H4").Value = "WEST"
$ws.Range("I4").Value = "SHWETA"
$ws.Range("J4").Formula = '=IF(B3>1000,"HIGH","LOW")'
$ws.Range("K4").Formula = '=IF(AND(B3>1000,H4="WEST"),"TOP PERFORMER","REGULAR")'
$ws.Range("L4").Formula = '=IF(OR(H4="WEST",H4="SOUTH"),"ZONE A","ZONE B")'
$ws.Range("M4").Formula = '=IFERROR(VLOOKUP(I4,F4:G8,2,FALSE),"NOT FOUND")'

# --- Row 5 -------------------------------------------------------------
$ws.Range("A5").Value = "ARJUN"
$ws.Range("B5").Value = 1600
$ws.Range("C5").Value = "HIGHEST SALE"
$ws.Range("D5").Formula = "=MAX(B5)"
$ws.Range("F5").Value = "PRIYA"
$ws.Range("G5").Value = 950
$ws.Range("H5").Value = "EAST"
$ws.Range("I5").Value = "SIMRAN"
$ws.Range("J5").Formula = '=IF(B4>1000,"HIGH","LOW")'
$ws.Range("K5").Formula = '=IF(AND(B4>1000,H5="WEST"),"TOP PERFORMER","REGULAR")'
$ws.Range("L5").Formula = '=IF(OR(H5="WEST",H5="SOUTH"),"ZONE A","ZONE B")'
$ws.Range("M5").Formula = '=IFERROR(VLOOKUP(I5,F5:G9,2,FALSE),"NOT FOUND")'

# --- Row 6 -------------------------------------------------------------
$ws.Range("A6").Value = "SIMRAN"
$ws.Range("B6").Value = 750
$ws.Range("C6").Value = "LOWEST SALE"
$ws.Range("D6").Formula = "=MIN(B3:B7)"
$ws.Range("F6").Value = "ARJUN"
$ws.Range("G6").Value = 1600
$ws.Range("H6").Value = "WEST"
$ws.Range("I6").Value = "RAHUL"
$ws.Range("J6").Formula = '=IF(B5>1000,"HIGH","LOW")'
$ws.Range("K6").Formula = '=IF(AND(B5>1000,H6="WEST"),"TOP PERFORMER","REGULAR")'
$ws.Range("L6").Formula = '=IF(OR(H6="WEST",H6="SOUTH"),"ZONE A","ZONE B")'
$ws.Range("M6").Formula = '=IFERROR(VLOOKUP(I6,F4:G10,2,FALSE),"NOT FOUND")'

# --- Row 7 -------------------------------------------------------------
$ws.Range("A7").Value = "NEHA"
$ws.Range("B7").Value = 1100
$ws.Range("C7").Value = "TOTAL PEOPLE"
$ws.Range("D7").Formula = "=COUNTA(B3:B8)"
$ws.Range("F7").Value = "SIMRAN"
$ws.Range("G7").Value = 750
$ws.Range("H7").Value = "SOUTH"
$ws.Range("I7").Value = "KAVYA"
$ws.Range("J7").Formula = '=IF(B6>1000,"HIGH","LOW")'
$ws.Range("K7").Formula = '=IF(AND(B6>1000,H7="WEST"),"TOP PERFORMER","REGULAR")'
$ws.Range("L7").Formula = '=IF(OR(H7="WEST",H7="SOUTH"),"ZONE A","ZONE B")'
$ws.Range("M7").Formula = '=IFERROR(VLOOKUP(I7,F7:G11,2,FALSE),"NOT FOUND")'

# --- Row 8 -------------------------------------------------------------
$ws.Range("B8").Value = "SNEHA"
$ws.Range("C8").Value = "HOW MANY NUMERIC VALUE "
$ws.Range("D8").Formula = "=COUNT(B3:B8)"
$ws.Range("D8").Borders.LineStyle = 1
$ws.Range("D8").Borders.Weight = 2
$ws.Range("F8").Value = "NEHA"
$ws.Range("G8").Value = 1100
$ws.Range("H8").Value = "NORTH"
$ws.Range("I8").Value = "NITYA"
$ws.Range("J8").Formula = '=IF(B7>1000,"HIGH","LOW")'
$ws.Range("K8").Formula = '=IF(AND(B7>1000,H8="WEST"),"TOP PERFORMER","REGULAR")'
$ws.Range("L8").Formula = '=IF(OR(H8="WEST",H8="SOUTH"),"ZONE A","ZONE B")'
$ws.Range("M8").Formula = '=IFERROR(VLOOKUP(I8,F8:G12,2,FALSE),"NOT FOUND")'

# --- Row 11 --------------------------------------------------------------
$ws.Range("A11").Value = "emp id "
$ws.Range("B11").Value = "name "

# --- Row 12 --------------------------------------------------------------
$ws.Range("A12").Value = 101
$ws.Range("B12").Value = "shweta"
$ws.Range("C12").Formula = "=VLOOKUP(101,A12:B15,2,FALSE)"
$ws.Range("D12").Value = "emp id"
$ws.Range("E12").Value = 101
$ws.Range("F12").Value = 102
$ws.Range("G12").Value = 103

# --- Row 13 --------------------------------------------------------------
$ws.Range("A13").Value = 102
$ws.Range("B13").Value = "shruti"
$ws.Range("C13").Formula = "=VLOOKUP(104,A12:B15,2,FALSE)"
$ws.Range("D13").Value = "name"
$ws.Range("E13").Value = "shweta"
$ws.Range("F13").Value = "rohit"
$ws.Range("G13").Value = "rahul"

# --- Row 14 --------------------------------------------------------------
$ws.Range("A14").Value = 103
$ws.Range("B14").Value = "sakshi"
$ws.Range("C14").Formula = "=XLOOKUP(103,A12:A15,B12:B15,2,FALSE)"

# --- Row 15 --------------------------------------------------------------
$ws.Range("A15").Value = 104
$ws.Range("B15").Value = "selu"
$ws.Range("C15").Formula = "=XLOOKUP(102,A12:A15,B12:B15,2,FALSE)"
$ws.Range("D15").Formula = "=HLOOKUP(102,E12:G13,2,FALSE)"
$ws.Range("E15").Formula = '=IFERROR(VLOOKUP(104,E12:G13,2,FALSE),"NOT FOUND")'
$ws.Range("I15").Formula = '=IF(VLOOKUP(I4,F4:G8,2,FALSE),"NOT FOUND")'

# Column widths to roughly match the authored layout.
$ws.Columns.Item(3).ColumnWidth = 29.14
$ws.Columns.Item(13).ColumnWidth = 20.14

# Selection / active tab: Sheet2 becomes the visible tab, cursor at K13.
$ws.Range("K13").Select()
